$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1828
$ws.Range("I40").Value = 1483.5
$ws.Range("J40").Value = 2241.4
$ws.Range("K40").Value = 1483.5
$ws.Range("L40").Value = 2241.4
$ws.Range("M40").Value = -1308.5
$ws.Range("N40").Value = -2591.4

$ws.Range("H86").Value = 1613
$ws.Range("I86").Value = 1867
$ws.Range("J86").Value = 1460.6
$ws.Range("K86").Value = 1867
$ws.Range("L86").Value = 1460.6
$ws.Range("M86").Value = -744
$ws.Range("N86").Value = -3706.6

$ws.Range("H89").Value = 1613
$ws.Range("I89").Value = 1867
$ws.Range("J89").Value = 1460.6
$ws.Range("K89").Value = 9335
$ws.Range("L89").Value = 7303
$ws.Range("M89").Value = -3719
$ws.Range("N89").Value = -18535

$ws.Range("H138").Value = 7939069.5
$ws.Range("I138").Value = 2803.0952
$ws.Range("J138").Value = 11907202
$ws.Range("K138").Value = 8409.285600000001
$ws.Range("L138").Value = 35721606
$ws.Range("M138").Value = -3269.285600000001
$ws.Range("N138").Value = -35731886

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2360.5652
$ws.Range("I61").Value = 1645.3684
$ws.Range("K61").Value = 1645.3684
$ws.Range("M61").Value = -1433.3684

$ws.Range("H122").Value = 2136.2693
$ws.Range("I122").Value = 2150.75
$ws.Range("K122").Value = 6452.25
$ws.Range("M122").Value = -4002.25

$ws.Range("H132").Value = 2397.9
$ws.Range("I132").Value = 1735.3334
$ws.Range("K132").Value = 5206.0002
$ws.Range("M132").Value = -2676.0002

$ws.Range("H136").Value = 2360.5652
$ws.Range("I136").Value = 1645.3684
$ws.Range("K136").Value = 4936.1052
$ws.Range("M136").Value = -2386.1052

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 15501.25
$ws.Range("I11").Value = 15501.25
$ws.Range("K11").Value = 15501.25
$ws.Range("M11").Value = -15361.25

$ws.Range("H16").Value = 350.5
$ws.Range("I16").Value = 374
$ws.Range("J16").Value = 280
$ws.Range("K16").Value = 374
$ws.Range("L16").Value = 280
$ws.Range("M16").Value = -87
$ws.Range("N16").Value = -854

$ws.Range("H31").Value = 1898.093
$ws.Range("I31").Value = 1198.05
$ws.Range("J31").Value = 2506.8262
$ws.Range("K31").Value = 1198.05
$ws.Range("L31").Value = 2506.8262
$ws.Range("M31").Value = -903.05
$ws.Range("N31").Value = -3096.8262

$ws.Range("H34").Value = 1898.093
$ws.Range("I34").Value = 1198.05
$ws.Range("J34").Value = 2506.8262
$ws.Range("K34").Value = 1198.05
$ws.Range("L34").Value = 2506.8262
$ws.Range("M34").Value = -996.05
$ws.Range("N34").Value = -2910.8262

$ws.Range("H39").Value = 3000
$ws.Range("I39").Value = 3000
$ws.Range("K39").Value = 3000
$ws.Range("M39").Value = -2609

$ws.Range("H49").Value = 3000
$ws.Range("I49").Value = 3000
$ws.Range("K49").Value = 3000
$ws.Range("M49").Value = -2818

$ws.Range("H93").Value = 16375
$ws.Range("I93").Value = 11166.667
$ws.Range("J93").Value = 32000
$ws.Range("K93").Value = 11166.667
$ws.Range("L93").Value = 32000
$ws.Range("M93").Value = -9294.666999999999
$ws.Range("N93").Value = -35744

$ws.Range("H99").Value = 8931129
$ws.Range("I99").Value = 12502740
$ws.Range("J99").Value = 2100
$ws.Range("K99").Value = 12502740
$ws.Range("L99").Value = 2100
$ws.Range("M99").Value = -12501242
$ws.Range("N99").Value = -5096

$ws.Range("H113").Value = 350.5
$ws.Range("I113").Value = 374
$ws.Range("J113").Value = 280
$ws.Range("K113").Value = 374
$ws.Range("L113").Value = 280
$ws.Range("M113").Value = 1796
$ws.Range("N113").Value = -4620

$ws.Range("H122").Value = 1122.8889
$ws.Range("I122").Value = 823
$ws.Range("J122").Value = 1362.8
$ws.Range("K122").Value = 2469
$ws.Range("L122").Value = 4088.4
$ws.Range("M122").Value = -19
$ws.Range("N122").Value = -8988.4

$ws.Range("H126").Value = 8931129
$ws.Range("I126").Value = 12502740
$ws.Range("J126").Value = 2100
$ws.Range("K126").Value = 37508220
$ws.Range("L126").Value = 6300
$ws.Range("M126").Value = -37505750
$ws.Range("N126").Value = -11240

$ws.Range("H132").Value = 3429.92
$ws.Range("I132").Value = 1922
$ws.Range("J132").Value = 5691.8
$ws.Range("K132").Value = 5766
$ws.Range("L132").Value = 17075.4
$ws.Range("M132").Value = -3236
$ws.Range("N132").Value = -22135.4

$ws.Range("H134").Value = 2271.111
$ws.Range("I134").Value = 704.7368
$ws.Range("J134").Value = 5991.25
$ws.Range("K134").Value = 2114.2104
$ws.Range("L134").Value = 17973.75
$ws.Range("M134").Value = 420.7896000000001
$ws.Range("N134").Value = -23043.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3842.6
$ws.Range("I68").Value = 4319.543
$ws.Range("J68").Value = 2173.3
$ws.Range("K68").Value = 12958.629
$ws.Range("L68").Value = 6519.900000000001
$ws.Range("M68").Value = -12147.629
$ws.Range("N68").Value = -8141.900000000001

$ws.Range("H71").Value = 3842.6
$ws.Range("I71").Value = 4319.543
$ws.Range("J71").Value = 2173.3
$ws.Range("K71").Value = 38875.887
$ws.Range("L71").Value = 19559.7
$ws.Range("M71").Value = -34819.887
$ws.Range("N71").Value = -27671.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H88").Value = 33695
$ws.Range("J88").Value = 33695
$ws.Range("L88").Value = 33695
$ws.Range("N88").Value = -34597

$ws.Range("H91").Value = 33695
$ws.Range("J91").Value = 33695
$ws.Range("L91").Value = 33695
$ws.Range("N91").Value = -36815

$ws.Range("H102").Value = 2405.7144
$ws.Range("I102").Value = 1904
$ws.Range("J102").Value = 2782
$ws.Range("K102").Value = 1904
$ws.Range("L102").Value = 2782
$ws.Range("M102").Value = -282
$ws.Range("N102").Value = -6026

$ws.Range("H132").Value = 2846.342
$ws.Range("I132").Value = 2203.9033
$ws.Range("K132").Value = 6611.7099
$ws.Range("M132").Value = -4081.7099

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3292083
$ws.Range("I100").Value = 10418765
$ws.Range("J100").Value = 2845.3845
$ws.Range("K100").Value = 10418765
$ws.Range("L100").Value = 2845.3845
$ws.Range("M100").Value = -10418224
$ws.Range("N100").Value = -3927.3845

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H50").Value = 15000
$ws.Range("J50").Value = 15000
$ws.Range("L50").Value = 15000
$ws.Range("N50").Value = -16262

$ws.Range("H57").Value = 39463.332
$ws.Range("J57").Value = 33695
$ws.Range("L57").Value = 33695
$ws.Range("N57").Value = -35203

$ws.Range("H132").Value = 11365560
$ws.Range("I132").Value = 14287114
$ws.Range("J132").Value = 3960.5557
$ws.Range("K132").Value = 42861342
$ws.Range("L132").Value = 11881.6671
$ws.Range("M132").Value = -42858812
$ws.Range("N132").Value = -16941.6671

$ws.Range("H136").Value = 16718402
$ws.Range("I136").Value = 22289592
$ws.Range("K136").Value = 66868776
$ws.Range("M136").Value = -66866226
